# working_hours.xlsx - "implemented handling of null values in the database"
#
# A new time-tracking entry (2014-07-30, 13:45 -> 18:00) is recorded in what
# used to be the blank placeholder row (row 167). A new blank placeholder row
# takes its place immediately below, and the three summary rows (sum [min],
# sum [h], sum [working weeks]) shift down by one row and their SUM() range
# grows to include the newly filled row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 167: the former row 167 (blank placeholder) and the
# three summary rows below it all shift down by one (167->168, 168->169,
# 169->170, 170->171).
$ws.Rows.Item(167).Insert()

# Populate the now-empty row 167 with the new working-hours entry.
$ws.Cells.Item(167, 1).Value = 2014
$ws.Cells.Item(167, 2).Value = 7
$ws.Cells.Item(167, 3).Value = 30
$ws.Cells.Item(167, 4).Value = 0.57291666666666663
$ws.Cells.Item(167, 5).Value = 0.75
$ws.Cells.Item(167, 6).Formula = "=(E167-D167)*24*60"
$ws.Cells.Item(167, 7).Formula = "=F167/60"

# The "sum [min]" total (now on row 169) needs to cover the newly added row.
$ws.Cells.Item(169, 6).Formula = "=SUM(F2:F167)"

# Match the author's final selection in the saved workbook.
$ws.Range("D168").Select()
